# Adding teachers to the Time Table workbook.
# Each course-code cell on the single "Time Table" sheet gets the
# assigned teacher initials appended in square brackets, e.g.
# "CS601  /  " -> "CS601[SSK]  /  " and paired-slot cells get both sides
# annotated, e.g. "CS693  /  CS691" -> "CS692[DC, SSK]  /  CS693[BDu, AP]".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "CS604A[DC]  /  CS604B[SLa]"
$ws.Range("B4").Value = "Free Period!"
$ws.Range("C4").Value = "Free Period!"
$ws.Range("E4").Value = "CS601[SSK]  /  "
$ws.Range("F4").Value = "CS602[DC]  /  "
$ws.Range("G4").Value = "CS605A[AH]  /  CS605B[SDe]"
$ws.Range("A6").Value = "CS692[DC, SSK]  /  CS693[BDu, AP]"
$ws.Range("B6").Value = "CS692[DC, SSK]  /  CS693[BDu, AP]"
$ws.Range("C6").Value = "CS692[DC, SSK]  /  CS693[BDu, AP]"
$ws.Range("D6").Value = "CS604A[DC]  /  CS604B[SLa]"
$ws.Range("E6").Value = "HU601[AnD]  /  "
$ws.Range("F6").Value = "CS603[BDu]  /  "
$ws.Range("G6").Value = "CS605A[AH]  /  CS605B[SDe]"
$ws.Range("A8").Value = "CS691[SSK, SDe]  /  CS692[DC, GY]"
$ws.Range("B8").Value = "CS691[SSK, SDe]  /  CS692[DC, GY]"
$ws.Range("C8").Value = "CS691[SSK, SDe]  /  CS692[DC, GY]"
$ws.Range("D8").Value = "CS601[SSK]  /  "
$ws.Range("E8").Value = "CS605A[AH]  /  CS605B[SDe]"
$ws.Range("F8").Value = "CS602[DC]  /  "
$ws.Range("G8").Value = "CS603[BDu]  /  "
$ws.Range("A10").Value = "Free Period!"
$ws.Range("B10").Value = "HU601[AnD]  /  "
$ws.Range("C10").Value = "CS603[BDu]  /  "
$ws.Range("E10").Value = "Free Period!"
$ws.Range("G10").Value = "Free Period!"
$ws.Range("A12").Value = "CS693[BDu, AP]  /  CS691[SSK, SDe]"
$ws.Range("B12").Value = "CS693[BDu, AP]  /  CS691[SSK, SDe]"
$ws.Range("C12").Value = "CS693[BDu, AP]  /  CS691[SSK, SDe]"
$ws.Range("D12").Value = "CS601[SSK]  /  "
$ws.Range("E12").Value = "CS604A[DC]  /  CS604B[SLa]"
$ws.Range("F12").Value = "CS602[DC]  /  "
$ws.Range("A16").Value = "IT605[AGh]  /  "
$ws.Range("B16").Value = "IT601A[AB]  /  IT601B[SL]"
$ws.Range("C16").Value = "IT603[ARC]  /  "
$ws.Range("D16").Value = "IT604A[AB]  /  IT604B[RG]"
$ws.Range("E16").Value = "IT692[ARC, KDa]  /  IT695[AGh, AB]"
$ws.Range("F16").Value = "IT692[ARC, KDa]  /  IT695[AGh, AB]"
$ws.Range("G16").Value = "IT692[ARC, KDa]  /  IT695[AGh, AB]"
$ws.Range("A18").Value = "IT605[AGh]  /  "
$ws.Range("B18").Value = "Free Period!"
$ws.Range("C18").Value = "Free Period!"
$ws.Range("E18").Value = "IT603[ARC]  /  "
$ws.Range("F18").Value = "Free Period!"
$ws.Range("G18").Value = "IT604A[AB]  /  IT604B[RG]"
$ws.Range("A20").Value = "IT693[SU, AD]  /  IT692[ARC, KDa]"
$ws.Range("B20").Value = "IT693[SU, AD]  /  IT692[ARC, KDa]"
$ws.Range("C20").Value = "IT693[SU, AD]  /  IT692[ARC, KDa]"
$ws.Range("E20").Value = "IT603[ARC]  /  "
$ws.Range("F20").Value = "IT602[SU]  /  "
$ws.Range("G20").Value = "HU601[SA]  /  "
$ws.Range("B22").Value = "HU685[ACh, SSR]  /  "
$ws.Range("C22").Value = "HU685[ACh, SSR]  /  "
$ws.Range("D22").Value = "HU685[ACh, SSR]  /  "
$ws.Range("E22").Value = "IT604A[AB]  /  IT604B[RG]"
$ws.Range("F22").Value = "IT602[SU]  /  "
$ws.Range("G22").Value = "IT601A[AB]  /  IT601B[SL]"
$ws.Range("A24").Value = "IT695[AGh, AB]  /  IT693[SU, AD]"
$ws.Range("B24").Value = "IT695[AGh, AB]  /  IT693[SU, AD]"
$ws.Range("C24").Value = "IT695[AGh, AB]  /  IT693[SU, AD]"
$ws.Range("D24").Value = "IT602[SU]  /  "
$ws.Range("E24").Value = "IT605[AGh]  /  "
$ws.Range("F24").Value = "HU601[SA]  /  "
$ws.Range("G24").Value = "IT601A[AB]  /  IT601B[SL]"
$ws.Range("A28").Value = "ECE602[TD]  /  "
$ws.Range("B28").Value = "ECE694[PP]  /  ECE695[SU]"
$ws.Range("C28").Value = "ECE694[PP]  /  ECE695[SU]"
$ws.Range("D28").Value = "ECE694[PP]  /  ECE695[SU]"
$ws.Range("E28").Value = "Free Period!"
$ws.Range("F28").Value = "ECE601[SSG]  /  "
$ws.Range("G28").Value = "ECE604[AnC]  /  "
$ws.Range("A30").Value = "ECE603A[PC]  /  ECE603B[JA]"
$ws.Range("B30").Value = "Free Period!"
$ws.Range("C30").Value = "Free Period!"
$ws.Range("D30").Value = "ECE602[TD]  /  "
$ws.Range("E30").Value = "HU685[ACh, AnC]  /  "
$ws.Range("F30").Value = "HU685[ACh, AnC]  /  "
$ws.Range("G30").Value = "HU685[ACh, AnC]  /  "
$ws.Range("A32").Value = "ECE603A[PC]  /  ECE603B[JA]"
$ws.Range("B32").Value = "ECE604[AnC]  /  "
$ws.Range("C32").Value = "ECE601[SSG]  /  "
$ws.Range("D32").Value = "Free Period!"
$ws.Range("E32").Value = "ECE605A[AD]  /  "
$ws.Range("F32").Value = "ECE602[TD]  /  "
$ws.Range("G32").Value = "HU601[TR]  /  "
$ws.Range("A34").Value = "ECE695[SU]  /  ECE692[SD]"
$ws.Range("B34").Value = "ECE695[SU]  /  ECE692[SD]"
$ws.Range("C34").Value = "ECE695[SU]  /  ECE692[SD]"
$ws.Range("D34").Value = "ECE604[AnC]  /  "
$ws.Range("E34").Value = "ECE605A[AD]  /  "
$ws.Range("F34").Value = "Free Period!"
$ws.Range("G34").Value = "HU601[TR]  /  "
$ws.Range("A36").Value = "ECE603A[PC]  /  ECE603B[JA]"
$ws.Range("B36").Value = "Free Period!"
$ws.Range("C36").Value = "ECE601[SSG]  /  "
$ws.Range("D36").Value = "ECE605A[AD]  /  "
$ws.Range("E36").Value = "ECE692[SD]  /  ECE694[PP]"
$ws.Range("F36").Value = "ECE692[SD]  /  ECE694[PP]"
$ws.Range("G36").Value = "ECE692[SD]  /  ECE694[PP]"
$ws.Range("A40").Value = "EE603(T)[SKB, SDG]  /  "
$ws.Range("B40").Value = "Free Period!"
$ws.Range("C40").Value = "Free Period!"
$ws.Range("D40").Value = "EE602[ABo]  /  "
$ws.Range("E40").Value = "EE603[SKB]  /  "
$ws.Range("F40").Value = "EE605A[IB]  /  EE605B[SD]"
$ws.Range("G40").Value = "Free Period!"
$ws.Range("A42").Value = "EE692[ABo, SDC]  /  EE691[PG, ASG]"
$ws.Range("B42").Value = "EE692[ABo, SDC]  /  EE691[PG, ASG]"
$ws.Range("C42").Value = "EE692[ABo, SDC]  /  EE691[PG, ASG]"
$ws.Range("D42").Value = "EE601[PG]  /  "
$ws.Range("E42").Value = "EE693[SDG, SMo]  /  EE692[ABo, SDC]"
$ws.Range("F42").Value = "EE693[SDG, SMo]  /  EE692[ABo, SDC]"
$ws.Range("G42").Value = "EE693[SDG, SMo]  /  EE692[ABo, SDC]"
$ws.Range("A44").Value = "EE602[ABo]  /  "
$ws.Range("B44").Value = "EE601[PG]  /  "
$ws.Range("C44").Value = "EE603[SKB]  /  "
$ws.Range("D44").Value = "EE604A[DC]  /  EE604B[GY]"
$ws.Range("E44").Value = "EE691[PG, ASG]  /  EE694B[AP, GY]"
$ws.Range("F44").Value = "EE691[PG, ASG]  /  EE694B[AP, GY]"
$ws.Range("G44").Value = "EE691[PG, ASG]  /  EE694B[AP, GY]"
$ws.Range("A46").Value = "EE604A[DC]  /  EE604B[GY]"
$ws.Range("B46").Value = "EE694B[GY, DC]  /  EE693[SDG, SMo]"
$ws.Range("C46").Value = "EE694B[GY, DC]  /  EE693[SDG, SMo]"
$ws.Range("D46").Value = "EE694B[GY, DC]  /  EE693[SDG, SMo]"
$ws.Range("E46").Value = "EE602(T)[ABo, SDC]  /  "
$ws.Range("F46").Value = "EE602[ABo]  /  "
$ws.Range("G46").Value = "EE605A[IB]  /  EE605B[SD]"
$ws.Range("A48").Value = "EE604A[DC]  /  EE604B[GY]"
$ws.Range("D48").Value = "EE605A[IB]  /  EE605B[JA]"
$ws.Range("E48").Value = "EE603[SKB]  /  "
$ws.Range("F48").Value = "EE601(T)[PG, ASG]  /  "
$ws.Range("G48").Value = "EE601[PG]  /  "
